# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the zh-cn and
# de-de report sheets to reflect a newly regenerated report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 03:59:34"
$wsZh.Range("H2").Value = "2016-03-19 04:00:16"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 03:59:42"
$wsDe.Range("H2").Value = "2016-03-19 04:00:32"
